$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 header text stays the same (re-set, same label text "HK_G_acc_LT")
$ws.Range("A1").Value = "HK_G_acc_LT"

# New set of accuracy values (rows 2-49), computed from the new HKG list run
$newValues = @(
    84.81481481481481,
    83.703703703703695,
    85.18518518518519,
    79.629629629629633,
    81.481481481481481,
    81.481481481481481,
    88.888888888888886,
    88.888888888888886,
    88.888888888888886,
    88.518518518518519,
    76.296296296296291,
    83.333333333333343,
    88.518518518518519,
    87.407407407407405,
    88.148148148148152,
    80,
    81.851851851851848,
    84.074074074074076,
    85.925925925925924,
    87.407407407407405,
    87.407407407407405,
    77.777777777777786,
    77.407407407407405,
    76.296296296296291,
    87.407407407407405,
    84.444444444444443,
    85.555555555555557,
    83.333333333333343,
    81.851851851851848,
    83.703703703703695,
    81.851851851851848,
    82.962962962962962,
    82.222222222222214,
    80,
    80.370370370370367,
    78.518518518518519,
    80.370370370370367,
    77.407407407407405,
    77.407407407407405,
    85.925925925925924,
    87.037037037037038,
    86.296296296296291,
    85.555555555555557,
    85.925925925925924,
    85.555555555555557,
    78.888888888888886,
    79.629629629629633,
    84.074074074074076
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
